$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 82 (pushes existing rows 82..185 down to 83..186)
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row with the new "Berenjena" market observation
$ws.Range("A82").Value = 5
$ws.Range("B82").Value = "Macroferia Regional de Talca"
$ws.Range("C82").Value = "Maule"
$ws.Range("D82").Value = 45128
$ws.Range("E82").Value = 7
$ws.Range("F82").Value = 100112001
$ws.Range("G82").Value = "Berenjena"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 150
$ws.Range("K82").Value = 8000
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = 8000
$ws.Range("N82").Value = "`$/caja 50 unidades"
$ws.Range("O82").Value = "Región de Arica y Parinacota"
$ws.Range("P82").Value = 160
$ws.Range("Q82").Value = 50
$ws.Range("R82").Value = "Hortaliza"
